$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage (Excel COM
# auto-coerces numeric-looking strings like "275.00" or "0.104" into real
# numbers on assignment, which would lose formatting/leading zeros/trailing
# zeros). Temporarily switching the NumberFormat to Text ("@") keeps the
# literal text, then the original Style is restored so no visible formatting
# change leaks into the workbook.
function Set-TextValue($rng, $val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '44.212.02'
$ws.Range("E2").Value = '  +1.15%  '
Set-TextValue $ws.Range("D3") '2.261.53'
$ws.Range("E3").Value = '  +2.69%  '
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue $ws.Range("D5") '99.39'
$ws.Range("E5").Value = '  +17.42%  '
Set-TextValue $ws.Range("D6") '275.00'
$ws.Range("E6").Value = '  +6.92%  '
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +6.22%  '
Set-TextValue $ws.Range("D10") '48.49'
$ws.Range("E10").Value = '  +6.79%  '
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("E12").Value = '  +13.30%  '
Set-TextValue $ws.Range("D13") '0.104'
$ws.Range("E13").Value = '  +0.40%  '
Set-TextValue $ws.Range("D14") '15.59'
$ws.Range("E14").Value = '  +8.40%  '
Set-TextValue $ws.Range("D15") '2.595.86'
$ws.Range("E15").Value = '  +2.50%  '
Set-TextValue $ws.Range("D16") '0.839'
$ws.Range("E16").Value = '  +6.69%  '
Set-TextValue $ws.Range("D17") '2.251.28'
$ws.Range("E17").Value = '  +3.65%  '
Set-TextValue $ws.Range("D18") '44.201.68'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("E19").Value = '  +3.51%  '
Set-TextValue $ws.Range("D20") '6.23'
$ws.Range("E20").Value = '  +5.20%  '
Set-TextValue $ws.Range("D21") '70.99'
$ws.Range("E21").Value = '  +1.48%  '
Set-TextValue $ws.Range("D22") '10.83'
$ws.Range("E22").Value = '  +20.55%  '
$ws.Range("E23").Value = '  -1.41%  '
Set-TextValue $ws.Range("D24") '235.34'
Set-TextValue $ws.Range("D26") '11.53'
$ws.Range("E26").Value = '  +8.21%  '
$ws.Range("E27").Value = '  +13.57%  '
Set-TextValue $ws.Range("D28") '40.22'
$ws.Range("E28").Value = '  +3.58%  '
Set-TextValue $ws.Range("D29") '3.39'
$ws.Range("E29").Value = '  -2.96%  '
Set-TextValue $ws.Range("D30") '2.29'
$ws.Range("E30").Value = '  +0.65%  '
Set-TextValue $ws.Range("D31") '173.55'
$ws.Range("E31").Value = '  -0.09%  '
Set-TextValue $ws.Range("D32") '0.0919'
$ws.Range("E32").Value = '  +6.52%  '
Set-TextValue $ws.Range("D33") '21.25'
$ws.Range("E33").Value = '  +4.05%  '
Set-TextValue $ws.Range("D34") '5.76'
$ws.Range("E34").Value = '  +8.23%  '
Set-TextValue $ws.Range("D35") '0.114'
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("E36").Value = '  +1.36%  '
Set-TextValue $ws.Range("D37") '0.0357'
$ws.Range("E37").Value = '  -0.97%  '
Set-TextValue $ws.Range("D38") '4.43'
$ws.Range("E38").Value = '  -1.28%  '
Set-TextValue $ws.Range("D39") '3.61'
$ws.Range("E39").Value = '  +25.62%  '
Set-TextValue $ws.Range("D40") '0.252'
$ws.Range("E40").Value = '  +26.68%  '
Set-TextValue $ws.Range("D41") '12.69'
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("E42").Value = '  +5.00%  '
Set-TextValue $ws.Range("D43") '62.84'
$ws.Range("E43").Value = '  -1.23%  '
Set-TextValue $ws.Range("D44") '5.50'
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("E45").Value = '  +5.73%  '
Set-TextValue $ws.Range("D46") '8.60'
$ws.Range("E46").Value = '  +3.05%  '
Set-TextValue $ws.Range("D47") '100.53'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E48").Value = '  +4.78%  '
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("E50").Value = '  +1.03%  '
Set-TextValue $ws.Range("D51") '2.478.31'
$ws.Range("E51").Value = '  +2.43%  '
